# Update the Registration sheet's sample user rows (emails/telephones bumped
# from the "...18"/"10000000[19-22]" batch to the "...22"/"10000000[31-34]" batch)
# and move the active-cell selection from E4 to E3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Registration")

# Row 2 (Ibraheem Abu)
$ws.Range("C2").Value = "ibu22@gmail.com"
$ws.Range("D2").Value = "'1000000031"

# Row 3 (Deepti Kharbanda)
$ws.Range("C3").Value = "deepti.kharbanda22@gmail.com"
$ws.Range("D3").Value = "'1000000032"

# Row 4 (Kartika Varma)
$ws.Range("C4").Value = "kartika.varma22@gmail.com"
$ws.Range("D4").Value = "'1000000033"

# Row 5 (Guddu Kharbanda)
$ws.Range("C5").Value = "kavita.kharbanda22@gmail.com"
$ws.Range("D5").Value = "'1000000034"

# Move the selection/active cell from E4 to E3 on the Registration sheet.
$ws.Activate() | Out-Null
$ws.Range("E3").Select() | Out-Null
